$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Challenge": the "Challenge Link" / "Credits" block in column F moves
# up by one row (the old answer cell "A2:A10." is removed), and the
# challenge-description text in F5 is reworded.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Challenge")

$xlPasteFormats = -4122

# Remove the two existing hyperlinks (they will be re-created lower down).
$ws.Hyperlinks.Delete()

# F5: reworded challenge description.
$ws.Range("F5").Value = "Provide a formula to sum last 3 non-zeros values in the range."

# F6 used to hold the old "A2:A10." answer text - it goes away completely.
$ws.Range("F6").Clear()

# F7 becomes the "Challenge Link" sub-heading (previously sat in F8).
$ws.Range("F4").Copy()
$ws.Range("F7").PasteSpecial($xlPasteFormats)
$ws.Range("F7").Value = "Challenge Link"

# F8 becomes the challenge-link hyperlink text (previously sat in F9).
$ws.Range("F9").Copy()
$ws.Range("F8").PasteSpecial($xlPasteFormats)
$ws.Range("F8").Value = "https://lnkd.in/dtw7Ni5u"

# F9 used to hold the challenge-link hyperlink text - now empty.
$ws.Range("F9").Clear()

# F10 becomes the "Credits" sub-heading (previously sat in F11).
$ws.Range("F4").Copy()
$ws.Range("F10").PasteSpecial($xlPasteFormats)
$ws.Range("F10").Value = "Credits"

# F11 becomes the credits hyperlink text (previously sat in F12).
$ws.Range("F12").Copy()
$ws.Range("F11").PasteSpecial($xlPasteFormats)
$ws.Range("F11").Value = "Excel BI"

# F12 used to hold the credits hyperlink text - now empty.
$ws.Range("F12").Clear()

$excel.CutCopyMode = $false

# Re-create the hyperlinks pointing at their new cells.
$ws.Hyperlinks.Add($ws.Range("F8"), "https://lnkd.in/dtw7Ni5u")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.linkedin.com/in/excelbi/")

# Hyperlinks.Add always re-stamps the "Hyperlink" cell style as a brand-new
# style record; flip it back onto the regular named "Hyperlink" style so the
# formatting matches the rest of the sheet (reuses the existing style slot).
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F11").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet "Power Query": record cell B36 as the active selection.
# ---------------------------------------------------------------------------
$pq = $wb.Worksheets.Item("Power Query")
$pq.Activate() | Out-Null
$pq.Range("B36").Select() | Out-Null

# Restore "Challenge" as the active/front-most sheet.
$ws.Activate() | Out-Null
